$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-128 down to 27-129.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new data record.
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44676
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = 100112026
$ws.Cells.Item(26, 7).Value = "Haba"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 73
$ws.Cells.Item(26, 11).Value = 18000
$ws.Cells.Item(26, 12).Value = 19000
$ws.Cells.Item(26, 13).Value = 18479
$ws.Cells.Item(26, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 739
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"

# Apply the date number format used by the rest of column D to the new D26 cell.
$ws.Cells.Item(26, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
